$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DoBlockAnalysis")

# Rename the shared-string text that F6 points at (LOOKUP_&_MAINTAIN_GLOBAL_VARS -> MAINTAIN_GLOBAL_VARS + USE_GLOBAL_VARS)
$ws.Range("F6").Value = "MAINTAIN_GLOBAL_VARS + USE_GLOBAL_VARS (Add additional vals to map var)"

# F27:F31 now reference the "CREATE_OUTPUT_PATTERN (called rule)" text instead of "RULE_INVOCATION_CONTROL (called rule)"
$createOutputPattern = "CREATE_OUTPUT_PATTERN" + [char]10 + " (called rule)"
$ws.Range("F27").Value = $createOutputPattern
$ws.Range("F28").Value = $createOutputPattern
$ws.Range("F29").Value = $createOutputPattern
$ws.Range("F30").Value = $createOutputPattern
$ws.Range("F31").Value = $createOutputPattern

# Update the sheet view: scroll/selection moved from E25/F27 to C10/F10
[void]$ws.Activate()
[void]$ws.Range("F10").Select()
